$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start time -> make it a live NOW() formula again
$ws.Range("B5").Formula = "=NOW()"

# Raised (By BSS -> To BSS)
$ws.Range("B8").Value = "To BSS"

# Module (NCC -> HLR)
$ws.Range("B10").Value = "HLR"

# Mail Subject
$ws.Range("B11").Value = "Linebarring Failed Orders | HLR"

# Issue Description
$ws.Range("B12").Value = "Linebarring Failed Orders ar HLR"

# Action Taken
$ws.Range("B13").Value = "Raised to HLR for WA"

# BMC TICKET ID
$ws.Range("B19").Value = "INC000000031410"

# Cosmetic: H29 drops its grey/italic-ish helper font back to the plain
# default font (matches the rest of the helper column, e.g. I7).
$ws.Range("H29").Font.Name = "Arial"
$ws.Range("H29").Font.Size = 10

# Data validation list updates (formula1 for a literal list needs to be a
# quoted, comma-separated string - mirrors how these were authored originally)
$ws.Range("B10").Validation.Delete()
$ws.Range("B10").Validation.Add(3, 1, 1, '"NCC,ERP,PRM,BSS,UMS,NMS,LMS,SND,Tibco,HLR"')

$ws.Range("B11").Validation.Delete()
$ws.Range("B11").Validation.Add(3, 1, 1, '"BSS SAFARICOM || FNF Order Failures,BSS SAFARICOM || Provisioning Failures from NCC,BSS Safaricom || Lifecyclesync Termination failure at HLR,Language update and Change SIM active failure cases,BSS Safaricom || ChangeSubscription failure at NCC,BSS Safari"')

$ws.Range("B12").Validation.Delete()
$ws.Range("B12").Validation.Add(3, 1, 1, '"ChangeSim order failure at NCC for WA,Provisioning Order Failure at NCC for WA,Provisioning order failures at SND for WA,Fnf Order Failure at NCC for WA,Connection Migration Failure at NCC for WA,Modifying existing customer profile  at BSS for WA,clear th"')

$ws.Range("B13").Validation.Delete()
$ws.Range("B13").Validation.Add(3, 1, 1, '"Raised to ERP for WA,Raised to NCC for WA,Raised to SND  for WA,Raised to Tibco  for WA,WA at BSS,Raised to HLR for WA"')

$ws.Range("B17").Validation.Delete()
$ws.Range("B17").Validation.Add(3, 1, 1, '"Asnake,Million,Abenezer,Seid,Akshay,Greejith,Arun,Anjali,Nishmita"')

# Sheet view: scroll back to the top and select A3
$ws.Range("A3").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
